$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, matching the style of the existing
# header cells (copy format from G1, the adjacent "sum" header).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Populate the Save column values for the existing data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
